$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (Group=uoa4, Year=2023, Currency=TRY, Value=20098221.34),
# mirroring the existing uoa3/TRY row's formatting.
$ws.Range("A5").Value = "uoa4"
$ws.Range("B5").Value = 2023
$ws.Range("C5").Value = "TRY"
$ws.Range("D5").Value = 20098221.34

# Copy the number formatting/style from D4 onto the new D5 cell.
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column D widens to fit the newly added values.
$ws.Columns.Item(4).AutoFit()

# Leave the selection where the user ended up after entering the row.
$ws.Range("F5").Select()
